$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 15:06:55"
$wsZhCn.Range("H2").Value = "2016-03-21 15:07:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 15:06:59"
$wsDeDe.Range("H2").Value = "2016-03-21 15:07:24"
